# DataCucu.xlsx edit
#
# Adds two new worksheets after the existing "Hoja1":
#   - "Hoja2"         : a blank sheet
#   - "InvalidFields"  : a copy of Hoja1's header row + its last (blank-template)
#                         data row, with a couple of fields overwritten with
#                         "invalid" test data, wrapped text on the email cell,
#                         and formatted as its own table ("Tabla13").
# It also updates the view state: Hoja1 is no longer the tab shown on open
# (InvalidFields is), and Hoja1's selection becomes the whole table range.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Create the new sheets, in tab order right after Hoja1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja2"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "InvalidFields"

# ---------------------------------------------------------------------------
# 2. Populate "InvalidFields" from Hoja1's header + template row, then edit
#    a few cells with bogus data (this is what produced the two new
#    "invalid fields" test rows referenced by the commit message).
# ---------------------------------------------------------------------------
$ws1.Range("A1:R1").Copy($ws3.Range("A1")) | Out-Null
$ws1.Range("A4:R4").Copy($ws3.Range("A2")) | Out-Null

$ws3.Range("M2").Value = "fgsdfgghfgh"
$ws3.Range("Q2").Value = "sgfhsghsfgh"

$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = "5465165"
$ws3.Range("C2").NumberFormat = "@"
$ws3.Range("C2").Value = "6344543"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "mailto:sadsada@hotmail.com") | Out-Null
$ws3.Range("A2").WrapText = $true
$ws3.Rows.Item(2).RowHeight = 28.8

# ---------------------------------------------------------------------------
# 3. Turn the populated range into its own table, mirroring Hoja1's table.
# ---------------------------------------------------------------------------
$lo = $ws3.ListObjects.Add(1, $ws3.Range("A1:R2"), $null, 1)
$lo.Name = "Tabla13"
$lo.TableStyle = "TableStyleMedium2"

# ---------------------------------------------------------------------------
# 4. View state: Hoja1 keeps the table selected (no longer the active tab),
#    InvalidFields becomes the active tab with P2 selected.
# ---------------------------------------------------------------------------
$ws1.Select()
$ws1.Range("A1:R4").Select()

$ws3.Select()
$ws3.Range("P2").Select()
